$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change per-row: D (4), L (12), M (13), N (14), O (15), P (16), S (19)
$cols = @(4, 12, 13, 14, 15, 16, 19)

# Snapshot the current values for every affected row/column BEFORE writing anything,
# since the update is a permutation of existing rows (row N's new content comes from
# a different row's old content).
$snapshot = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row mapping (content of target row becomes what source row used to hold)
$mapping = @{
    2  = 10
    3  = 2
    4  = 3
    5  = 6
    6  = 7
    7  = 9
    8  = 4
    9  = 8
    10 = 5
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $sourceVals[$c]
    }
}
